$d = $word.ActiveDocument

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$rsq = [char]0x2019   # right single quotation mark (curly apostrophe)

# Replaces the contents of the (unique) paragraph whose visible text equals
# $oldText with the full paragraph markup supplied in $xml. Using InsertXML
# on the paragraph's own Range (rather than Range.Text = "...") keeps
# sibling empty runs (<w:r/>) and paragraph formatting intact, since
# Range.Text assignment collapses/merges the paragraph's runs.
function Set-ParagraphXml($oldText, $xml) {
    foreach ($p in $d.Paragraphs) {
        $t = $p.Range.Text
        if ($t.Length -gt 0) {
            $trimmed = $t.Substring(0, $t.Length - 1)
        } else {
            $trimmed = $t
        }
        if ($trimmed -eq $oldText) {
            $p.Range.InsertXML($xml) | Out-Null
            return $true
        }
    }
    Write-Host "WARNING: paragraph not found for [$oldText]"
    return $false
}

$newTitle = "Play Joker" + $rsq + "s Five Free - Classic Slot Game with Joker Wild Symbol"

# --- Title heading (Heading1) ---
Set-ParagraphXml "Play Joker's Five Slot for Free - Review & Pros/Cons" `
    "<w:p $wns><w:pPr><w:pStyle w:val=`"Heading1`"/></w:pPr><w:r><w:t>$newTitle</w:t></w:r></w:p>" | Out-Null

# --- "What we like" bullet list ---
Set-ParagraphXml "Classic slot game with potential for interesting wins" `
    "<w:p $wns><w:pPr><w:pStyle w:val=`"ListBullet`"/><w:spacing w:line=`"240`" w:lineRule=`"auto`"/><w:ind w:left=`"720`"/></w:pPr><w:r/><w:r><w:t>Classic slot game with 5 paylines</w:t></w:r></w:p>" | Out-Null

Set-ParagraphXml "Cute and simple graphics in a comfortable layout" `
    "<w:p $wns><w:pPr><w:pStyle w:val=`"ListBullet`"/><w:spacing w:line=`"240`" w:lineRule=`"auto`"/><w:ind w:left=`"720`"/></w:pPr><w:r/><w:r><w:t>Joker wild symbol with 5x multiplier</w:t></w:r></w:p>" | Out-Null

Set-ParagraphXml "Medium volatility and a good RTP" `
    "<w:p $wns><w:pPr><w:pStyle w:val=`"ListBullet`"/><w:spacing w:line=`"240`" w:lineRule=`"auto`"/><w:ind w:left=`"720`"/></w:pPr><w:r/><w:r><w:t>Cute and simple graphics</w:t></w:r></w:p>" | Out-Null

Set-ParagraphXml "Joker wild symbol that multiplies winning combinations" `
    "<w:p $wns><w:pPr><w:pStyle w:val=`"ListBullet`"/><w:spacing w:line=`"240`" w:lineRule=`"auto`"/><w:ind w:left=`"720`"/></w:pPr><w:r/><w:r><w:t>Medium volatility with a good RTP</w:t></w:r></w:p>" | Out-Null

# --- "What we don't like" bullet list ---
Set-ParagraphXml "The gameplay is relatively traditional" `
    "<w:p $wns><w:pPr><w:pStyle w:val=`"ListBullet`"/><w:spacing w:line=`"240`" w:lineRule=`"auto`"/><w:ind w:left=`"720`"/></w:pPr><w:r/><w:r><w:t>Big wins may come with longer intervals</w:t></w:r></w:p>" | Out-Null

# --- Bold title repeated near the end ---
Set-ParagraphXml "Play Joker's Five Slot for Free - Review & Pros/Cons" `
    "<w:p $wns><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>$newTitle</w:t></w:r></w:p>" | Out-Null

# --- Italic summary near the end ---
$newItalic = "Play Joker" + $rsq + "s Five for free and enjoy a classic slot game experience with the Joker wild symbol."
Set-ParagraphXml "Explore the features of Joker's Five, a classic online slot game from Synot Games. Read our review and play for free to experience the Joker wild symbol." `
    "<w:p $wns><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>$newItalic</w:t></w:r></w:p>" | Out-Null
